$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new row (row 5) of data to Sheet1.
$ws.Range("A5").Value = "444DDD444"

# Format B5 as Text first so the date-looking string "2025-10-23" is
# stored as a literal string rather than being auto-converted into a
# date serial number.
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2025-10-23"

$ws.Range("C5").Value = "Velachery"
$ws.Range("D5").Value = "DDD"
